# adding profits to tables
# Adds a third "M_PL" block (columns R:Y) mirroring the existing
# "M_%cit" (B:I) / "M_ETR" (J:Q) blocks: a merged header label in row 1,
# the same 8 group-column headers in row 2, and new profit figures in
# rows 4-8 (row 7 only has data through column U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcCols = @("J","K","L","M","N","O","P","Q")
$dstCols = @("R","S","T","U","V","W","X","Y")

# --- Row 1: merged header cell R1:Y1 = "M_PL", styled like the other
#     two header blocks (bold + border + centered). Merge the (still
#     blank/default-styled) range first, THEN paste the formatting —
#     merging an already-bordered range makes Excel redistribute the
#     border per merged-cell position instead of leaving every cell
#     with the same style id. ---
$ws.Range("R1:Y1").Merge()
$ws.Range("J1:Q1").Copy()
$ws.Range("R1:Y1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("R1").Value = "M_PL"

# --- Row 2: same 8 group headers repeated under the new block,
#     copying formatting from the J2:Q2 block first. ---
$ws.Range("J2:Q2").Copy()
$ws.Range("R2:Y2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $ws.Range("$($dstCols[$i])2").Value = $ws.Range("$($srcCols[$i])2").Value()
}

# --- Rows 4-8: new profit values for the M_PL block. ---
$row4 = @(958366954868, 958617846484, 956175459402, 956426351018, 1008209699708, 1008209699708, 1008209699708, 1008209699708)
$row5 = @(3140810, 3140810, 3140810, 3140810, 3140810, 3140810, 3140810, 3140810)
$row6 = @(2557154421, 16062039693, 2124037026, 17056843395, 21277927825, 21277927825, 21277927825, 21277927825)
$row7 = @(11464456004, 868911225, 1750040641, 868911225)
$row8 = @(41567527900, 47033047618, 41286747272, 47534570308, 49527932043, 49527932043, 49527932043, 49527932043)

for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Range("$($dstCols[$i])4").Value = $row4[$i]
}
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Range("$($dstCols[$i])5").Value = $row5[$i]
}
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Range("$($dstCols[$i])6").Value = $row6[$i]
}
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Range("$($dstCols[$i])7").Value = $row7[$i]
}
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Range("$($dstCols[$i])8").Value = $row8[$i]
}
